$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ingredients "selegiline" and "sulpiride" are no longer part of this
# list (it now only tracks antidepressant ingredients, not the
# non-opioid/antipsychotic extras that had crept in). Remove those two
# rows; Excel shifts every following row up automatically and drops the
# now-unused shared strings on save.
$ws.Range("A45").EntireRow.Delete() | Out-Null   # "selegiline"
$ws.Range("A47").EntireRow.Delete() | Out-Null   # "sulpiride" (shifted up after the first delete)

# Re-apply the existing custom sort (the list keeps its tail section,
# rows 40 onward, sorted alphabetically) so the sort range bookkeeping
# reflects the new, smaller extent of the sheet.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A40:A56")) | Out-Null
$ws.Sort.SetRange($ws.Range("A40:A56"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Restore the selection to where the user was last working.
$ws.Range("D48").Select()
